$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)

function Replace-WholeText($tr, [string]$old, [string]$new) {
    $full = $tr.Text
    $idx = $full.IndexOf($old)
    if ($idx -lt 0) {
        throw "Could not find text: $old"
    }
    $sub = $tr.Characters($idx + 1, $old.Length)
    $sub.Text = $new
}

# --- Shape "Rectangle 2": square -> rectangle wording in the 4 worked examples ---
$sh2 = $s.Shapes.Item("Rectangle 2")
$tr2 = $sh2.TextFrame.TextRange

Replace-WholeText $tr2 "# Calculate the area of a square." "# Calculate the area of a rectangle."
Replace-WholeText $tr2 ' square with sides 3m by 5m has area { 3 * 5 }m squared")' ' rectangle with sides 3m by 5m has area { 3 * 5 }m squared")'
Replace-WholeText $tr2 "# Calculate the area of a different square." "# Calculate the area of a different rectangle."
Replace-WholeText $tr2 ' square with sides 2m by 6m has area { 2 * 6 }m squared")' ' rectangle with sides 2m by 6m has area { 2 * 6 }m squared")'
Replace-WholeText $tr2 "# Calculate the area of yet another different square." "# Calculate the area of yet another different rectangle."
Replace-WholeText $tr2 ' square with sides 5m by 7m has area { 5 * 7 }m squared")' ' rectangle with sides 5m by 7m has area { 5 * 7 }m squared")'
Replace-WholeText $tr2 "# Calculate the area of a final square" "# Calculate the area of a final rectangle"
Replace-WholeText $tr2 ' square with sides 4m by 4m has area { 4 * 4 }m squared")' ' rectangle with sides 4m by 4m has area { 4 * 4 }m squared")'

# Text re-wrapped onto extra lines now that "rectangle" is longer than "square";
# resize the autofit textbox to match.
$sh2.Height = 3162404 / 12700

# --- Shape "Rectangle 3": the print_area_square function + its calls ---
$sh3 = $s.Shapes.Item("Rectangle 3")
$tr3 = $sh3.TextFrame.TextRange

# Re-typing "def print_area_square" merged the old "def" / " " runs into one "def " run.
$full3 = $tr3.Text
$idxDef = $full3.IndexOf("def")
$rDef = $tr3.Characters($idxDef + 1, 3)
$rSpace = $tr3.Characters($idxDef + 1 + 3, 1)
$rSpace.Text = "def "
$rDef.Text = ""

Replace-WholeText $tr3 "print_area_square" "print_area_rectangle"
Replace-WholeText $tr3 ' square with sides {a}m by {b}m has area {a * b}m squared")' ' rectangle with sides {a}m by {b}m has area {a * b}m squared")'
Replace-WholeText $tr3 "print_area_square" "print_area_rectangle"
Replace-WholeText $tr3 "print_area_square" "print_area_rectangle"
Replace-WholeText $tr3 "print_area_square" "print_area_rectangle"
Replace-WholeText $tr3 "print_area_square" "print_area_rectangle"

$sh3.Height = 1708160 / 12700
